$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2 and 3): column B (CONTENT) text now matches
# the question TYPE instead of the generic "fill in the blank" wording ---
$ws.Range("B2").Value2 = "Chọn một đáp án đúng"
$ws.Range("B3").Value2 = "Chọn đáp án đúng"

# --- Add two new sample rows for the "tf" (true/false) and "fb" (fill the
# blank) question types ---
$ws.Range("A4").Value2 = "Câu hỏi 3"
$ws.Range("B4").Value2 = "Chọn đáp án đúng"
$ws.Range("C4").Value2 = "easy"
$ws.Range("D4").Value2 = "tf"
$ws.Range("E4").Value2 = "0;"
$ws.Range("F4").Value2 = "1;"

$ws.Range("A5").Value2 = "Câu hỏi 4"
$ws.Range("B5").Value2 = "Điền từ vào chỗ trống"
$ws.Range("C5").Value2 = "hard"
$ws.Range("D5").Value2 = "fb"
$ws.Range("E5").Value2 = "a;an"
$ws.Range("F5").Value2 = "are;is"

# --- Column widths (widened to fit the longer CONTENT/CORRECT_OPTION/
# WRONG_OPTION text; values chosen to land as close as possible to the
# target stored widths 22.7109375 / 18.28515625 / 18.140625 once the host
# snaps the assigned character width to its internal column-width grid) ---
$ws.Columns.Item(2).ColumnWidth = 21.833333333333336
$ws.Columns.Item(5).ColumnWidth = 17.5
$ws.Columns.Item(6).ColumnWidth = 17.333333333333336

# --- Selection ends on H5, matching the author's last cursor position ---
$ws.Range("H5").Select() | Out-Null
